$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/b46f630a61e4ca49f359104ae5b3caff8584e07d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/61cfc6ec573a6d5974f85cc7e1504a898362f577/e2e/b.md."

# --- Overview sheet: update the "b.md" row (row 3) with the new handoff status/date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-01-03 04:48:40"

# --- zh-cn sheet: update the "b.md" row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-01-03 04:48:29"
$wsZhCn.Range("R3").Value = $errorDetail
$wsZhCn.Columns.Item(18).ColumnWidth = 40

# --- de-de sheet: update the "b.md" row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-01-03 04:48:40"
$wsDeDe.Range("R3").Value = $errorDetail
$wsDeDe.Columns.Item(18).ColumnWidth = 40
